$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new row ("ACI") was inserted at row 38, pushing ACIFORMULA..RECKITTBEN down by one.
$ws.Rows.Item(38).Insert()

# Copy the label formatting (bold, centered, bordered) from a neighboring label cell
# onto the newly inserted row so it matches the rest of column A.
$ws.Cells.Item(37,1).Copy()
$ws.Cells.Item(38,1).PasteSpecial(-4122)

# Refresh every data row (2-45) with the recalculated statistics, including the
# newly-inserted "ACI" row and the appended "JMISMDL" row at the end.
$ws.Cells.Item(2,1).Value = "ISLAMIBANK"
$ws.Cells.Item(2,2).Value = 0.007583692224203503
$ws.Cells.Item(2,3).Value = 0.4391619553311629
$ws.Cells.Item(2,4).Value = 0.194572513446317
$ws.Cells.Item(2,5).Value = 0.04514955934143219
$ws.Cells.Item(3,1).Value = "RAKCERAMIC"
$ws.Cells.Item(3,2).Value = 0.007234200582994451
$ws.Cells.Item(3,3).Value = 0.4188033614631229
$ws.Cells.Item(3,4).Value = 0.4293144845826624
$ws.Cells.Item(3,5).Value = 0.0880238353129781
$ws.Cells.Item(4,1).Value = "BATBC"
$ws.Cells.Item(4,2).Value = -0.001341162784693873
$ws.Cells.Item(4,3).Value = 0.5871709323758949
$ws.Cells.Item(4,4).Value = 0.210594985777164
$ws.Cells.Item(4,5).Value = 0.09844525387137125
$ws.Cells.Item(5,1).Value = "RENATA"
$ws.Cells.Item(5,2).Value = -0.002361204303685866
$ws.Cells.Item(5,3).Value = 0.4382668741649042
$ws.Cells.Item(5,4).Value = 0.1519476336520654
$ws.Cells.Item(5,5).Value = 0.0555209587766058
$ws.Cells.Item(6,1).Value = "MARICO"
$ws.Cells.Item(6,2).Value = 0.005915072252129519
$ws.Cells.Item(6,3).Value = 0.5446635761522398
$ws.Cells.Item(6,4).Value = 0.1067136247020446
$ws.Cells.Item(6,5).Value = 0.02654926617877797
$ws.Cells.Item(7,1).Value = "SINGERBD"
$ws.Cells.Item(7,2).Value = -0.008259500407661293
$ws.Cells.Item(7,3).Value = 0.5929733520400338
$ws.Cells.Item(7,4).Value = 0.1722810247887226
$ws.Cells.Item(7,5).Value = 0.08014817881042202
$ws.Cells.Item(8,1).Value = "LINDEBD"
$ws.Cells.Item(8,2).Value = 0.003209620601093824
$ws.Cells.Item(8,3).Value = 0.6973500738677669
$ws.Cells.Item(8,4).Value = 0.2781467624323112
$ws.Cells.Item(8,5).Value = 0.08361292333238486
$ws.Cells.Item(9,1).Value = "SUMITPOWER"
$ws.Cells.Item(9,2).Value = -0.009949035924624398
$ws.Cells.Item(9,3).Value = 0.5752935888060694
$ws.Cells.Item(9,4).Value = 0.2244549116666649
$ws.Cells.Item(9,5).Value = 0.1156103251459202
$ws.Cells.Item(10,1).Value = "BERGERPBL"
$ws.Cells.Item(10,2).Value = 0.0006950183342122699
$ws.Cells.Item(10,3).Value = 0.4531813946218111
$ws.Cells.Item(10,4).Value = 0.1874646958355251
$ws.Cells.Item(10,5).Value = 0.0213303314342743
$ws.Cells.Item(11,1).Value = "BXPHARMA"
$ws.Cells.Item(11,2).Value = -0.007587514766900601
$ws.Cells.Item(11,3).Value = 0.6232552650927289
$ws.Cells.Item(11,4).Value = 0.3604053787687483
$ws.Cells.Item(11,5).Value = 0.1592061054918643
$ws.Cells.Item(12,1).Value = "SQURPHARMA"
$ws.Cells.Item(12,2).Value = -0.00100913236631307
$ws.Cells.Item(12,3).Value = 0.655809616366313
$ws.Cells.Item(12,4).Value = 0.1516993362048366
$ws.Cells.Item(12,5).Value = 0.04071489019504188
$ws.Cells.Item(13,1).Value = "SAPORTL"
$ws.Cells.Item(13,2).Value = [double]"-2.781304057549564e-17"
$ws.Cells.Item(13,3).Value = 0.4747563278230189
$ws.Cells.Item(13,4).Value = 0.5906645243792101
$ws.Cells.Item(13,5).Value = 0.08986426341005095
$ws.Cells.Item(14,1).Value = "BSCCL"
$ws.Cells.Item(14,2).Value = 0.01127495055271219
$ws.Cells.Item(14,3).Value = 0.5402638429872133
$ws.Cells.Item(14,4).Value = 0.3510998139591944
$ws.Cells.Item(14,5).Value = 0.07368820104220086
$ws.Cells.Item(15,1).Value = "GP"
$ws.Cells.Item(15,2).Value = -0.009439810767149778
$ws.Cells.Item(15,3).Value = 0.5483915725761019
$ws.Cells.Item(15,4).Value = 0.1614687073359332
$ws.Cells.Item(15,5).Value = 0.1067477523407605
$ws.Cells.Item(16,1).Value = "SQUARETEXT"
$ws.Cells.Item(16,2).Value = 0.02161092246603014
$ws.Cells.Item(16,3).Value = 0.4842984192449983
$ws.Cells.Item(16,4).Value = 0.3592992085137503
$ws.Cells.Item(16,5).Value = 0.1274182744398169
$ws.Cells.Item(17,1).Value = "OLYMPIC"
$ws.Cells.Item(17,2).Value = -0.005201524107832577
$ws.Cells.Item(17,3).Value = 0.4234608676886672
$ws.Cells.Item(17,4).Value = 0.422842654482525
$ws.Cells.Item(17,5).Value = 0.1509119980366686
$ws.Cells.Item(18,1).Value = "BEXIMCO"
$ws.Cells.Item(18,2).Value = 0.008414960840092458
$ws.Cells.Item(18,3).Value = 0.5724283454106756
$ws.Cells.Item(18,4).Value = 0.3993797653925039
$ws.Cells.Item(18,5).Value = 0.1553701867004263
$ws.Cells.Item(19,1).Value = "HFL"
$ws.Cells.Item(19,2).Value = 0.0002551752009745702
$ws.Cells.Item(19,3).Value = 0.5366806416913728
$ws.Cells.Item(19,4).Value = 0.7277247958205673
$ws.Cells.Item(19,5).Value = 0.1524096567952328
$ws.Cells.Item(20,1).Value = "UPGDCL"
$ws.Cells.Item(20,2).Value = -0.007511533223774774
$ws.Cells.Item(20,3).Value = 0.715160658591476
$ws.Cells.Item(20,4).Value = 0.1816790117080007
$ws.Cells.Item(20,5).Value = 0.08385397402405664
$ws.Cells.Item(21,1).Value = "DELTALIFE"
$ws.Cells.Item(21,2).Value = 0.01782946512821283
$ws.Cells.Item(21,3).Value = 0.3630246513360748
$ws.Cells.Item(21,4).Value = 0.609327814031273
$ws.Cells.Item(21,5).Value = 0.1650132386094825
$ws.Cells.Item(22,1).Value = "BENGALWTL"
$ws.Cells.Item(22,2).Value = 0.01597194860349844
$ws.Cells.Item(22,3).Value = 0.4713055996315534
$ws.Cells.Item(22,4).Value = 0.6339683267960111
$ws.Cells.Item(22,5).Value = 0.1656049803862032
$ws.Cells.Item(23,1).Value = "EHL"
$ws.Cells.Item(23,2).Value = 0.03444764892126757
$ws.Cells.Item(23,3).Value = 0.6058994412300485
$ws.Cells.Item(23,4).Value = 0.7865582176225234
$ws.Cells.Item(23,5).Value = 0.324055638460285
$ws.Cells.Item(24,1).Value = "IBNSINA"
$ws.Cells.Item(24,2).Value = 0.007069059272421766
$ws.Cells.Item(24,3).Value = 0.4697966559174029
$ws.Cells.Item(24,4).Value = 0.2665632567513026
$ws.Cells.Item(24,5).Value = 0.05220968130771937
$ws.Cells.Item(25,1).Value = "WALTONHIL"
$ws.Cells.Item(25,2).Value = -0.01187637180980473
$ws.Cells.Item(25,3).Value = 0.5844133271923974
$ws.Cells.Item(25,4).Value = 0.2115586032731094
$ws.Cells.Item(25,5).Value = 0.09860025916436817
$ws.Cells.Item(26,1).Value = "UNIQUEHRL"
$ws.Cells.Item(26,2).Value = 0.03382567128801861
$ws.Cells.Item(26,3).Value = 0.5827544092068262
$ws.Cells.Item(26,4).Value = 0.6191915697149668
$ws.Cells.Item(26,5).Value = 0.163807649345543
$ws.Cells.Item(27,1).Value = "UNILEVERCL"
$ws.Cells.Item(27,2).Value = -0.01304571144913306
$ws.Cells.Item(27,3).Value = 0.2591811750769609
$ws.Cells.Item(27,4).Value = 0.5586423727616112
$ws.Cells.Item(27,5).Value = 0.09649148159035213
$ws.Cells.Item(28,1).Value = "KDSALTD"
$ws.Cells.Item(28,2).Value = 0.02124270175216034
$ws.Cells.Item(28,3).Value = 0.4336508899641016
$ws.Cells.Item(28,4).Value = 0.5075241217314312
$ws.Cells.Item(28,5).Value = 0.1671720362408985
$ws.Cells.Item(29,1).Value = "JHRML"
$ws.Cells.Item(29,2).Value = 0.09789403484022299
$ws.Cells.Item(29,3).Value = 0.3301858231836048
$ws.Cells.Item(29,4).Value = 0.9152162180837783
$ws.Cells.Item(29,5).Value = 0.2147727238845759
$ws.Cells.Item(30,1).Value = "ADNTEL"
$ws.Cells.Item(30,2).Value = 0.04189473134768912
$ws.Cells.Item(30,3).Value = 0.3827840727738596
$ws.Cells.Item(30,4).Value = 0.7000644644727796
$ws.Cells.Item(30,5).Value = 0.3883541633954137
$ws.Cells.Item(31,1).Value = "ITC"
$ws.Cells.Item(31,2).Value = 0.006191535800342367
$ws.Cells.Item(31,3).Value = 0.4988583452383596
$ws.Cells.Item(31,4).Value = 0.5177021720913639
$ws.Cells.Item(31,5).Value = 0.08938425441998503
$ws.Cells.Item(32,1).Value = "SIMTEX"
$ws.Cells.Item(32,2).Value = 0.01133251815313998
$ws.Cells.Item(32,3).Value = 0.4834132539420783
$ws.Cells.Item(32,4).Value = 0.5794821386053124
$ws.Cells.Item(32,5).Value = 0.1398929925479702
$ws.Cells.Item(33,1).Value = "APEXFOODS"
$ws.Cells.Item(33,2).Value = 0.05444105491848872
$ws.Cells.Item(33,3).Value = 0.4808440865345769
$ws.Cells.Item(33,4).Value = 0.7919269011369534
$ws.Cells.Item(33,5).Value = 0.3486936155836917
$ws.Cells.Item(34,1).Value = "FORTUNE"
$ws.Cells.Item(34,2).Value = 0.0256555146860071
$ws.Cells.Item(34,3).Value = 0.3909758398169211
$ws.Cells.Item(34,4).Value = 0.5231595754958333
$ws.Cells.Item(34,5).Value = 0.2484870924471906
$ws.Cells.Item(35,1).Value = "INDEXAGRO"
$ws.Cells.Item(35,2).Value = -0.002978096695202381
$ws.Cells.Item(35,3).Value = 0.4199274932301824
$ws.Cells.Item(35,4).Value = 0.581718050752609
$ws.Cells.Item(35,5).Value = 0.09412142235800772
$ws.Cells.Item(36,1).Value = "AMANFEED"
$ws.Cells.Item(36,2).Value = -0.01769132438505762
$ws.Cells.Item(36,3).Value = 0.5122829099016731
$ws.Cells.Item(36,4).Value = 0.5332546105291154
$ws.Cells.Item(36,5).Value = 0.2882314238924919
$ws.Cells.Item(37,1).Value = "SKICL"
$ws.Cells.Item(37,2).Value = 0.0839262632317533
$ws.Cells.Item(37,3).Value = 0.3385178070268697
$ws.Cells.Item(37,4).Value = 0.7556863948345267
$ws.Cells.Item(37,5).Value = 0.2111474088239285
$ws.Cells.Item(38,1).Value = "ACI"
$ws.Cells.Item(38,2).Value = -0.000400851658327197
$ws.Cells.Item(38,3).Value = 0.5625146278608716
$ws.Cells.Item(38,4).Value = 0.2721795536009472
$ws.Cells.Item(38,5).Value = 0.06698261552086486
$ws.Cells.Item(39,1).Value = "ACIFORMULA"
$ws.Cells.Item(39,2).Value = 0.00625771586139604
$ws.Cells.Item(39,3).Value = 0.5031060981873214
$ws.Cells.Item(39,4).Value = 0.4408664370492797
$ws.Cells.Item(39,5).Value = 0.06154815193748095
$ws.Cells.Item(40,1).Value = "EXIMBANK"
$ws.Cells.Item(40,2).Value = -0.002138827387020643
$ws.Cells.Item(40,3).Value = 0.5026779733311746
$ws.Cells.Item(40,4).Value = 0.2310046605304765
$ws.Cells.Item(40,5).Value = 0.09243796023272054
$ws.Cells.Item(41,1).Value = "IFILISLMF1"
$ws.Cells.Item(41,2).Value = -0.002761369937808542
$ws.Cells.Item(41,3).Value = 0.4739599818267394
$ws.Cells.Item(41,4).Value = 0.3125402711377219
$ws.Cells.Item(41,5).Value = 0.04024591790824537
$ws.Cells.Item(42,1).Value = "AIBL1STIMF"
$ws.Cells.Item(42,2).Value = -0.0115163095334896
$ws.Cells.Item(42,3).Value = 0.4249728655669395
$ws.Cells.Item(42,4).Value = 0.1946011773023534
$ws.Cells.Item(42,5).Value = 0.09416072109146312
$ws.Cells.Item(43,1).Value = "EXIM1STMF"
$ws.Cells.Item(43,2).Value = -0.01109665800925568
$ws.Cells.Item(43,3).Value = 0.6364454014058804
$ws.Cells.Item(43,4).Value = 0.3575631832035486
$ws.Cells.Item(43,5).Value = 0.1389273199405706
$ws.Cells.Item(44,1).Value = "RECKITTBEN"
$ws.Cells.Item(44,2).Value = 0.005815903023764328
$ws.Cells.Item(44,3).Value = 0.5902138311706845
$ws.Cells.Item(44,4).Value = 0.2202964518151829
$ws.Cells.Item(44,5).Value = 0.07836322939846488
$ws.Cells.Item(45,1).Value = "JMISMDL"
$ws.Cells.Item(45,2).Value = -0.011043144800542
$ws.Cells.Item(45,3).Value = 0.4848253665546084
$ws.Cells.Item(45,4).Value = 0.5538092716125815
$ws.Cells.Item(45,5).Value = 0.1611839577230219

# The appended JMISMDL row also needs the standard label formatting applied.
$ws.Cells.Item(44,1).Copy()
$ws.Cells.Item(45,1).PasteSpecial(-4122)
$ws.Cells.Item(45,1).Value = "JMISMDL"
